$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.233.32'
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('D3').Value = '2.279.64'
$ws.Range('E3').Value = '  -0.75%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.22'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.62'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.05%  '
$ws.Range('E7').Value = '  -1.14%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.592'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.68'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.38%  '
$ws.Range('E11').Value = '  -1.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.17'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.85%  '
$ws.Range('E13').Value = '  +1.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.972'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.00'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.77%  '
$ws.Range('D16').Value = '2.625.10'
$ws.Range('E16').Value = '  -0.75%  '
$ws.Range('D17').Value = '2.285.59'
$ws.Range('E17').Value = '  -0.48%  '
$ws.Range('D18').Value = '42.185.26'
$ws.Range('E18').Value = '  -0.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.23'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.48%  '
$ws.Range('E20').Value = '  -1.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.37'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.75'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.92%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.49'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '262.68'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.78%  '
$ws.Range('E25').Value = '  -4.44%  '
$ws.Range('E26').Value = '  +0.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.61'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.97%  '
$ws.Range('E28').Value = '  -1.94%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.87'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +12.99%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.31'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.84%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.87'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '163.77'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0858'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.96%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.129'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.89%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.60'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.20%  '
$ws.Range('E36').Value = '  -5.60%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.46'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.38%  '
$ws.Range('E38').Value = '  -4.08%  '
$ws.Range('E39').Value = '  -0.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.65'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.72%  '
$ws.Range('E41').Value = '  +4.80%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '97.09'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.44%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '68.80'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.31%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.225'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.999'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.15%  '
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('D47').Value = '1.699.72'
$ws.Range('E47').Value = '  +6.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '79.04'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.16%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '109.70'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.50%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.17'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.74%  '
$ws.Range('B51').Value = 'FraxShare'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.61'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.93%  '
